$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row -> (A, B, D, E) new values
$updates = @{
    2  = @(29, 7, 3, "В")
    3  = @(28, 7, 3, "В")
    4  = @(27, 6, 3, "СВ")
    5  = @(26, 6, 3, "СВ")
    6  = @(25, 6, 3, "СВ")
    7  = @(24, 6, 3, "СВ")
    8  = @(23, 7, 2, "С")
    9  = @(22, 7, 2, "С")
    10 = @(21, 7, 2, "С")
    11 = @(20, 7, 2, "С")
}

foreach ($row in $updates.Keys) {
    $vals = $updates[$row]
    $ws.Cells.Item($row, 1).Value = $vals[0]
    $ws.Cells.Item($row, 2).Value = $vals[1]
    $ws.Cells.Item($row, 4).Value = $vals[2]
    $ws.Cells.Item($row, 5).Value = $vals[3]
}
